$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 52 -------------------------------------------------------------
# event_id / fecha are stored as text (not auto-converted to number/date),
# matching the source feed's inlineStr cells. Force text via NumberFormat
# "@" for the assignment, then drop the leftover number-format style so the
# cell keeps its default (unstyled) appearance.
$ws.Range("A52").NumberFormat = "@"
$ws.Range("A52").Value = "14584706"
$ws.Range("A52").ClearFormats()

$ws.Range("B52").NumberFormat = "@"
$ws.Range("B52").Value = "2025-09-01"
$ws.Range("B52").ClearFormats()

$ws.Range("C52").Value = "Kilian Feldbausch"
$ws.Range("D52").Value = "Lorenzo Carboni"
$ws.Range("E52").Value = "Gana Lorenzo Carboni"
$ws.Range("F52").Value = 2

# --- Row 53 -------------------------------------------------------------
$ws.Range("A53").NumberFormat = "@"
$ws.Range("A53").Value = "14581498"
$ws.Range("A53").ClearFormats()

$ws.Range("B53").NumberFormat = "@"
$ws.Range("B53").Value = "2025-09-01"
$ws.Range("B53").ClearFormats()

$ws.Range("C53").Value = "Corentin Denolly"
$ws.Range("D53").Value = "Sandro Kopp"
$ws.Range("E53").Value = "Gana Corentin Denolly"
$ws.Range("F53").Value = 3.4
